# Append two new arrival rows (row 26 and row 27) to the "Main Data" sheet,
# matching the data already present for earlier flights on Jan 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Wizz Air A321 from London (LTN), tail (G-WUKI)
$ws.Cells.Item(26, 1).Value2 = 25
$ws.Cells.Item(26, 2).Value2 = "Sunday, Jan 15"
$ws.Cells.Item(26, 3).Value2 = "9:35 AM"
$ws.Cells.Item(26, 4).Value2 = "W92066"
$ws.Cells.Item(26, 5).Value2 = "London"
$ws.Cells.Item(26, 6).Value2 = "(LTN)"
$ws.Cells.Item(26, 7).Value2 = "Wizz Air "
$ws.Cells.Item(26, 8).Value2 = "A321"
$ws.Cells.Item(26, 9).Value2 = "(G-WUKI)"
$ws.Cells.Item(26, 10).Value2 = "9:14 AM"
$ws.Cells.Item(26, 12).Value2 = "0 hours, -21 minutes"

# Row 27: Wizz Air A321 from Eindhoven (EIN), tail (HA-LTC)
$ws.Cells.Item(27, 1).Value2 = 26
$ws.Cells.Item(27, 2).Value2 = "Sunday, Jan 15"
$ws.Cells.Item(27, 3).Value2 = "12:05 PM"
$ws.Cells.Item(27, 4).Value2 = "W62091"
$ws.Cells.Item(27, 5).Value2 = "Eindhoven"
$ws.Cells.Item(27, 6).Value2 = "(EIN)"
$ws.Cells.Item(27, 7).Value2 = "Wizz Air "
$ws.Cells.Item(27, 8).Value2 = "A321"
$ws.Cells.Item(27, 9).Value2 = "(HA-LTC)"
$ws.Cells.Item(27, 10).Value2 = "11:54 AM"
$ws.Cells.Item(27, 12).Value2 = "0 hours, -11 minutes"
